$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 41-92 were reordered into 5 contiguous blocks (grouped by sample-batch date
# prefix in column A), while the data within each block kept its relative order:
#   20230328 (25 rows), 20221229 (5), 20230227 (2), 20230322 (10), 20230215 (10)
# Columns B..K hold per-condition mean_abundance figures for each sample.
$data = @(
    @("20230328_BC49", 0, 0.4734714817043348, 0.2367758489022886, 0.03461341248054982, -0.06627049737934505, 0.405222155164785, -0.03138990012000637, 0, 0.2128625484694311, -0.01611671625890349),
    @("20230328_BC34", 0.003242027769600459, -0.3055183303786462, -0.08975896611411276, 0.01563619740350616, -0.02504774728019, 0, -0.05847624972055927, 0, 0.1498593426501968, -0.0642650581136136),
    @("20230328_BC36", 0.0002340432487577688, 0.3815807064938851, -0.1978249993904991, 0.0005668567284404699, -0.00037790906870502, 0, -0.00173854100843925, 0, 0.005013971239921531, -0.07387079684215367),
    @("20230328_BC32", 0.09304865387277343, 0.06807742493905188, 0.01904300450032795, 0.04326801355565556, 0.03099233014914076, 0.001312401093446732, -0.1089681114114057, 0, 0.3349538862559152, 0.02242734089608086),
    @("20230328_BC50", 0, 0.599542856884803, 0.3602220426005968, 0.0412572431457547, -0.03214344749897289, 0.6299842377202385, -0.03469041747507113, 0, 0.1473185935807473, 0.01009839842068523),
    @("20230328_BC38", 0.0005007129434112753, 0.3582835190100669, 0.2105141703694496, 0.03553814276429519, -0.01444460447922911, 0.0003768351015841419, -0.07234279212903365, 0, 0.3275854741911098, 0.0141445630425926),
    @("20230328_BC56", 0, 0.3994443391903282, 0.2002143629713586, 0.003053960828695146, -0.00005751579568346572, 0.000309812224699517, -0.0006178647303781026, 0, 0.0163723632144271, -0.000227546315732848),
    @("20230328_BC40", 0.0001259887988451976, -0.3983897609429662, -0.1995439304146367, 0.0008430762712905687, 0.001615504135889955, 0.00141516323001678, -0.000565475756808927, 0, 0.003977040273632854, -0.07419869639469365),
    @("20230328_BC46", 0.0001806694327884162, 0.4297223318867707, -0.1840462697536613, 0.00242700722454102, 0.0005343661507276878, 0.03777994047901675, -0.003962228858374002, 0, 0.002076829775157875, -0.0739414252118215),
    @("20230328_BC39", 0.1154463489834404, 0.001096478202762636, 0.009573111895565417, 0.04367250421021832, -0.001253003118877752, 0.0003437291010828508, -0.1509944902386528, 0, 0.3204712196277874, 0.009784075835095878),
    @("20230328_BC55", 0, -0.3957196886854581, -0.1992174067054509, 0.0002355827429344371, 0.0005309146145305087, 0.00035476969375929, -0.003242859427732992, 0, 0.001614584127191518, -0.07409770364425095),
    @("20230328_BC37", 0.0002018700335050141, 0.6444063348965056, 0.1923901183158563, 0.04362778695403654, 0.001692343325596965, 0.8023676594734058, -0.0009046411744440336, 0, 0.0570425870118039, -0.04801326472917215),
    @("20230328_BC52", 0, 0.2613338714237741, -0.1675843979287207, 0.006814577840159208, 0.001590482118348704, 0.07472850348686792, -0.01101803859812299, 0, 0.03162223332087603, -0.07257731563583986),
    @("20230328_BC53", 0, 0.3480025442825344, -0.1934535125318542, 0.003188068091590521, 0.002904600243471131, 0.01270200713055006, -0.003751061464202092, 0, 0.01319354670463498, -0.07409659481279253),
    @("20230328_BC35", 0.04961142045900987, 0.00521900605442657, -0.03030455920622626, 0.03428844518603862, 0.08254731748980279, 0.00208613048860985, -0.08617920515860394, 0, 0.3453831580341765, -0.01538845856351616),
    @("20230328_BC44", 0.0008607191975030398, 0.3982977124970158, -0.1923549108087214, 0.005256361543136298, -0.001824722585110354, 0.0009745631715955816, -0.04099566325857641, 0, 0.02449763176891431, -0.06987807766289345),
    @("20230328_BC42", 0.004335913778937655, 0.002088621066367305, -0.08759632199366532, 0.07202300690066152, 0.1219805963843645, 0.001308322086717353, -0.01505293590549694, 0, 0.424794654996318, -0.04986328666780647),
    @("20230328_BC43", 0.04360489125221342, 0.004035841735210067, 0.06173171734080975, 0.05205639607142023, 0.02769907979125615, 0.001959263201339136, -0.03461377842243839, 0, 0.2889378793093548, 0.03909494746881032),
    @("20230328_BC47", 0.002925565762898168, 0.5250743423531566, -0.1075265477378467, 0.02661626415265027, 0.003702147873254721, 0.1726275241549894, -0.02993349114000885, 0, 0.08781089186231646, -0.0683812632474522),
    @("20230328_BC41", 0.001148637417649912, 0.7213107669487043, 0.4084712694191846, 0.03389207231252128, -0.02181285141755355, 0.5752444255527747, -0.03121150527249909, 0, 0.1183917037325343, -0.0001174233531782534),
    @("20230328_BC33", 0.0002943184400109342, 0.3809792779422981, 0.06035340739613961, 0.006146375233451682, -0.0002233724440588878, 0.001645535127311687, -0.00196282619532578, 0, 0.07074068026922395, -0.04913896965747732),
    @("20230328_BC51", 0.0001660872515703498, 0.4061654044813175, -0.1970753044955014, 0.0006493603258158104, 0.0007005257599454669, 0.005474723952584086, -0.002461028807268576, 0, 0.003385590539910413, -0.07391170681891351),
    @("20230328_BC54", 0, 0.2497498330109085, -0.108064479692432, 0.01151681891871075, 0.00009428185745657239, 0.2216295086061052, 0, 0, 0.01098680432613265, -0.07307955902222539),
    @("20230328_BC48", 0.02053145260785224, 0.01013993207877523, -0.003052275376347325, 0.03898369620588998, 0.006087165013871916, 0.003711011290748371, -0.2119666022357595, 0, 0.2665100180545621, -0.0024980311445477),
    @("20230328_BC45", 0.0002361391667761761, -0.08864603235418382, -0.1743894817776503, 0.01154451775673332, 0.003115222304964215, 0.05802032297637293, -0.00141167699072312, 0, 0.09422268086451825, -0.0734192933009388),
    @("20221229_BC05", 0.0001436887058596133, 0.006809120875666342, 0.02514169713097846, 0.02956382631082029, 0.08362274592650137, 0.004259878554775546, -0.04841314364517252, 0.0004310661175788398, 0.06337060012279559, 0.005530279562489603),
    @("20221229_BC02", 0, 0.9073523165051274, 0.1645150926535315, 0.02142998904424739, 0.0006608411074830628, 0.4337580668915035, -0.0001251219659195118, 0, 0.006675376245480656, -0.05360303970627268),
    @("20221229_BC07", 0.0003515090021242112, 0.02070657309915077, 0.02277345523051209, 0.02938209225055256, 0.08574466944034755, 0.01418887026202229, -0.003444643568369848, 0, 0.03257072340170612, 0.004713708404295343),
    @("20221229_BC04", 0.0001083674969436144, 0.6774008519921054, 0.333744987315867, 0.0232582869225356, 0.002173060150725628, 0.4479313446485098, -0.00008132125344261773, 0, 0.01612437489357143, -0.001871308369608463),
    @("20221229_BC06", 0.004732347989563723, 0.05728468825514191, 0.07704071163639019, 0.03991601464315513, 0.02272316899607529, 0.04000247883859348, -0.1755030642603472, 0, 0.455754082002866, 0.04103531872845018),
    @("20230227_BC07", 0.0007646094702336669, 0, 0.1604386497538714, 0.05357564501321219, -0.06466713108977493, 0, -0.1668373262989449, 0, 0.537458327250021, 0.05012886020817157),
    @("20230227_BC08", 0.0002462494751640889, -0.3931005107260633, -0.1884938488984786, 0.01173239164674763, -0.0005400670468844388, 0, -0.04130275202309252, 0, 0.07494473944002632, -0.05921010167765459),
    @("20230322_BC72", 0.0005957349376948467, 0.3266341149358714, -0.1967936468330777, 0.0006879170818814007, 0.0008564011504021497, 0, 0, 0, 0.006222209588455295, -0.07358326332373645),
    @("20230322_BC77", 0.00470366355490682, 0.3224087571479483, 0.1200129272720805, 0.005685359355822619, 0.002312787483989106, 0.0004999254641340411, -0.01117159670343016, 0, 0.0233251146188738, -0.01660256829037488),
    @("20230322_BC89", 0.0004047578036891061, 0.2979716544032026, -0.1435404476265331, 0.02425899510943187, 0.003094762268407197, 0, -0.00180634262519586, 0, 0.283731721222013, -0.05466572447187669),
    @("20230322_BC73", 0.0003774950209381719, 0.001360310974638831, 0.1532091752974166, 0.0711303332508099, -0.05671926542422598, 0.0004255124605938259, -0.03677417963493045, 0, 0.5587231085931922, 0.05967478774968606),
    @("20230322_BC74", 0.0002050786870548268, 0.7866420849621016, 0.3030157815907856, 0.04160919023925777, 0.007159346532333234, 0.7260773203249894, -0.007631127057487481, 0, 0.09878366251778299, -0.02563966161033644),
    @("20230322_BC75", 0, -0.01585101907006612, -0.1964023562735602, 0.003106612292561146, 0.0003329371532924941, 0.002142976721455944, -0.002777550374411077, 0, 0.02410709313842841, -0.07032824245064531),
    @("20230322_BC76", 0.0003750451586889561, 0.4332033026983651, -0.1787125979886822, 0.004453606566888041, -0.0002742354187623928, 0.04226783572174848, 0, 0, 0.02784455675744881, -0.07394310855212725),
    @("20230322_BC96", 0.0005777345143148615, 0.4148077457482883, -0.1871570575900215, 0.005353380527920511, -0.00009493152118558392, 0.02018333166142914, -0.0181434029254432, 0, 0.0279664888000852, -0.06988218092605747),
    @("20230322_BC79", 0.001595805326623886, 0.0003208562489406237, 0.1066224305048903, 0.055263236796379, -0.02675364813148197, 0, -0.06750474379283376, 0, 0.3242803035519092, 0.05135114727595195),
    @("20230322_BC71", 0.1436741946998149, -0.0002044413135630802, 0.001986327905227178, 0.003388561787408536, 0.07570906212642746, 0, 0, 0, 0.01513288341148863, 0.0004847669301326959),
    @("20230215_BC05", 0.0007516627235206429, 0.002431375943846321, 0.167290913396781, 0.0584486412172475, -0.1904573785617034, 0.002362092658060304, -0.06661655346863167, 0, 0.4234863355985217, 0.07057382742844523),
    @("20230215_BC02", 0.001781252832546947, -0.3965682382373086, -0.1956095547619678, 0.001534040700747923, 0.003103319152292319, 0.001432428234580979, -0.0007649669725093023, 0, 0.006008401483458241, -0.07391959897034961),
    @("20230215_BC03", 0.00202358367917119, 0.5641045363008861, 0.2252554387906617, 0.04043734058626419, -0.03461410909283938, 0.431269310956221, -0.03460242536108083, 0, 0.2249651484635211, 0.002800715477050127),
    @("20230215_BC10", 0, -0.3955749949190014, -0.1976109848077313, 0.0005796516545710772, -0.000003327171235515791, 0.00168635246023749, 0, 0, 0.004063030136549752, -0.0739292558085282),
    @("20230215_BC07", 0.001720038275192313, 0.06363508834757217, 0.03388231272059633, 0.02308410362769227, -0.0001256424861763877, 0.4551593729400856, -0.0002360041662053635, 0, 0.005891497190374688, -0.0570442080358546),
    @("20230215_BC11", 0.001957558868288111, 0.002226147616308051, 0.1051989447989606, 0.04035558673171258, -0.1239008701492063, 0.003850106634696882, -0.1803340268049415, 0, 0.3263301820015905, 0.06307213783363309),
    @("20230215_BC01", 0.03611284804440207, 0.02637653436105057, 0.008534447192460697, 0.03520654611182204, 0.005778279507964328, 0.01430495978994945, -0.2868010864950036, 0, 0.2945679672334278, 0.01312985809421407),
    @("20230215_BC04", 0.0006713675252128478, 0.3997218760819036, -0.1978455978394082, 0.001063917224430386, 0.0003267453224374919, 0.001205448382074343, -0.001000380974857434, 0, 0.003234855303132942, -0.07398031752322884),
    @("20230215_BC06", 0.0002644620240375748, 0.6668522803353798, 0.3173842983266412, 0.02297244027443103, -0.0002194665714951181, 0.4647165114701988, -0.0003966930360563622, 0, 0.003016142899326173, -0.006285478067598878),
    @("20230215_BC08", 0.000874677648181117, 0.6607669630310343, 0.1525008840123738, 0.02526157837597555, 0.0001219798136466278, 0.4807035510890998, -0.02118757719097042, 0, 0.02201613105549413, -0.05360255875489292)
)

$startRow = 41
$r = $startRow
foreach ($row in $data) {
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}